$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 371.66666
$ws.Range("I28").Value = 371.66666
$ws.Range("K28").Value = 371.66666
$ws.Range("M28").Value = 113.33334

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 118
$ws.Range("I41").Value = 127
$ws.Range("K41").Value = 127
$ws.Range("M41").Value = 313

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 13411.583
$ws.Range("I74").Value = 11988.4
$ws.Range("J74").Value = 14428.143
$ws.Range("K74").Value = 11988.4
$ws.Range("L74").Value = 14428.143
$ws.Range("M74").Value = -11052.4
$ws.Range("N74").Value = -16300.143

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3196
$ws.Range("I76").Value = 3248
$ws.Range("K76").Value = 3248
$ws.Range("M76").Value = -2933

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 13411.583
$ws.Range("I77").Value = 11988.4
$ws.Range("J77").Value = 14428.143
$ws.Range("K77").Value = 59942
$ws.Range("L77").Value = 72140.715
$ws.Range("M77").Value = -55262
$ws.Range("N77").Value = -81500.715

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3196
$ws.Range("I79").Value = 3248
$ws.Range("K79").Value = 3248
$ws.Range("M79").Value = -2156

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3250.4
$ws.Range("I116").Value = 3100.6667
$ws.Range("J116").Value = 3475
$ws.Range("K116").Value = 3100.6667
$ws.Range("L116").Value = 3475
$ws.Range("M116").Value = 341.3332999999998
$ws.Range("N116").Value = -10359

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1725.1034
$ws.Range("I137").Value = 1564
$ws.Range("K137").Value = 4692
$ws.Range("M137").Value = -2142

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3199.8
$ws.Range("I88").Value = 2999.6667
$ws.Range("K88").Value = 2999.6667
$ws.Range("M88").Value = -2593.6667

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3199.8
$ws.Range("I91").Value = 2999.6667
$ws.Range("K91").Value = 2999.6667
$ws.Range("M91").Value = -1595.6667

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1493.5217
$ws.Range("I132").Value = 1472.65
$ws.Range("K132").Value = 4417.950000000001
$ws.Range("M132").Value = -1887.950000000001

# BSM row 8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 25252
$ws.Range("J8").Value = 49500
$ws.Range("L8").Value = 49500
$ws.Range("N8").Value = -49780

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2641.4
$ws.Range("I105").Value = 2511.5789
$ws.Range("J105").Value = 3052.5
$ws.Range("K105").Value = 2511.5789
$ws.Range("L105").Value = 3052.5
$ws.Range("M105").Value = -764.5789
$ws.Range("N105").Value = -6546.5

# BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 99914.664
$ws.Range("J132").Value = 99914.664
$ws.Range("L132").Value = 99914.664
$ws.Range("N132").Value = -110034.664

# CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4500.091
$ws.Range("I6").Value = 2250.5
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 2250.5
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = -2137.5
$ws.Range("N6").Value = -5226

# CRP row 28
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 58724.75
$ws.Range("J28").Value = 58724.75
$ws.Range("L28").Value = 58724.75
$ws.Range("N28").Value = -59214.75

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12643.471
$ws.Range("I58").Value = 9323.333000000001
$ws.Range("K58").Value = 9323.333000000001
$ws.Range("M58").Value = -9120.333000000001

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4857.143
$ws.Range("J62").Value = 4375
$ws.Range("L62").Value = 4375
$ws.Range("N62").Value = -5623

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4857.143
$ws.Range("J65").Value = 4375
$ws.Range("L65").Value = 21875
$ws.Range("N65").Value = -28115

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4377.533
$ws.Range("I99").Value = 4174
$ws.Range("K99").Value = 4174
$ws.Range("M99").Value = -2676

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4377.533
$ws.Range("I126").Value = 4174
$ws.Range("K126").Value = 12522
$ws.Range("M126").Value = -10052

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 12643.471
$ws.Range("I136").Value = 9323.333000000001
$ws.Range("K136").Value = 27969.999
$ws.Range("M136").Value = -25419.999

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 174478.8
$ws.Range("J141").Value = 198207.31
$ws.Range("L141").Value = 198207.31
$ws.Range("N141").Value = -208567.31

# GSM row 31
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 15000
$ws.Range("I31").Value = 15000
$ws.Range("K31").Value = 15000
$ws.Range("M31").Value = -14708

# GSM row 37
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 15000
$ws.Range("I37").Value = 15000
$ws.Range("K37").Value = 15000
$ws.Range("M37").Value = -14723

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 18000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 18000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -18540

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 18000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 18000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 18000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -19872

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 921.381
$ws.Range("J97").Value = 1194.6666
$ws.Range("L97").Value = 1194.6666
$ws.Range("N97").Value = -2186.6666

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 523.5625
$ws.Range("J16").Value = 1164.5
$ws.Range("L16").Value = 1164.5
$ws.Range("N16").Value = -1504.5

# LTW row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 10250
$ws.Range("I32").Value = 10250
$ws.Range("K32").Value = 10250
$ws.Range("M32").Value = -9933

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7236.75
$ws.Range("I46").Value = 7236.75
$ws.Range("K46").Value = 7236.75
$ws.Range("M46").Value = -7048.75

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 8210.888999999999
$ws.Range("J100").Value = 8333.333000000001
$ws.Range("L100").Value = 8333.333000000001
$ws.Range("N100").Value = -9415.333000000001

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2237.25
$ws.Range("I2").Value = 966
$ws.Range("K2").Value = 966
$ws.Range("M2").Value = -854

# WVR row 39
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 22222
$ws.Range("I39").Value = 22222
$ws.Range("K39").Value = 22222
$ws.Range("M39").Value = -21809

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1866.6666
$ws.Range("I96").Value = 1800
$ws.Range("J96").Value = 1900
$ws.Range("K96").Value = 1800
$ws.Range("L96").Value = 1900
$ws.Range("M96").Value = -427
$ws.Range("N96").Value = -4646

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5641.5
$ws.Range("I132").Value = 4641.56
$ws.Range("K132").Value = 13924.68
$ws.Range("M132").Value = -11394.68
